$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.73885133333334
$ws.Range("H2").Value = 107.216554
$ws.Range("I2").Value = 0.01949729408921566
$ws.Range("J2").Value = 0.01949729408921566
$ws.Range("M2").Value = 0.7489546666666667
$ws.Range("N2").Value = 2.246864
$ws.Range("O2").Value = 0.05220789806691288
$ws.Range("P2").Value = 0.05220789806691287
$ws.Range("Q2").Value = 26.76677948740623
$ws.Range("R2").Value = 240.901015386656
$ws.Range("S2").Value = 0.001017912742390394
$ws.Range("T2").Value = 0.001017912742390394
$ws.Range("G3").Value = 35.73885133333334
$ws.Range("H3").Value = 107.216554
$ws.Range("I3").Value = 0.01949729408921566
$ws.Range("J3").Value = 0.01949729408921566
$ws.Range("O3").Value = 0.1982273102638064
$ws.Range("P3").Value = 0.1982273102638064
$ws.Range("Q3").Value = 101.6303451905416
$ws.Range("R3").Value = 914.6731067148741
$ws.Range("S3").Value = 0.003864896164727632
$ws.Range("T3").Value = 0.003864896164727631
$ws.Range("G4").Value = 35.73885133333334
$ws.Range("H4").Value = 107.216554
$ws.Range("I4").Value = 0.01949729408921566
$ws.Range("J4").Value = 0.01949729408921566
$ws.Range("M4").Value = 10.337765
$ws.Range("N4").Value = 31.013295
$ws.Range("O4").Value = 0.7206216949842531
$ws.Range("P4").Value = 0.720621694984253
$ws.Range("Q4").Value = 369.4598464539367
$ws.Range("R4").Value = 3325.13861808543
$ws.Range("S4").Value = 0.01405017311417705
$ws.Range("T4").Value = 0.01405017311417704
$ws.Range("G5").Value = 35.73885133333334
$ws.Range("H5").Value = 107.216554
$ws.Range("I5").Value = 0.01949729408921566
$ws.Range("J5").Value = 0.01949729408921566
$ws.Range("M5").Value = 0.4152066666666667
$ws.Range("N5").Value = 1.24562
$ws.Range("O5").Value = 0.02894309668502767
$ws.Range("P5").Value = 0.02894309668502767
$ws.Range("Q5").Value = 14.83900933260889
$ws.Range("R5").Value = 133.55108399348
$ws.Range("S5").Value = 0.0005643120679205875
$ws.Range("T5").Value = 0.0005643120679205872
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.9215900675332435
$ws.Range("J6").Value = 0.9215900675332435
$ws.Range("M6").Value = 0.7489546666666667
$ws.Range("N6").Value = 2.246864
$ws.Range("O6").Value = 0.05220789806691288
$ws.Range("P6").Value = 0.05220789806691287
$ws.Range("Q6").Value = 1265.201109578098
$ws.Range("R6").Value = 11386.80998620288
$ws.Range("S6").Value = 0.04811428030525493
$ws.Range("T6").Value = 0.04811428030525492
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.9215900675332435
$ws.Range("J7").Value = 0.9215900675332435
$ws.Range("O7").Value = 0.1982273102638064
$ws.Range("P7").Value = 0.1982273102638064
$ws.Range("Q7").Value = 4803.821302535724
$ws.Range("R7").Value = 43234.39172282152
$ws.Range("S7").Value = 0.1826843202529546
$ws.Range("T7").Value = 0.1826843202529546
$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.9215900675332435
$ws.Range("J8").Value = 0.9215900675332435
$ws.Range("M8").Value = 10.337765
$ws.Range("N8").Value = 31.013295
$ws.Range("O8").Value = 0.7206216949842531
$ws.Range("P8").Value = 0.720621694984253
$ws.Range("Q8").Value = 17463.47586933293
$ws.Range("R8").Value = 157171.2828239964
$ws.Range("S8").Value = 0.6641177965464582
$ws.Range("T8").Value = 0.6641177965464581
$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.9215900675332435
$ws.Range("J9").Value = 0.9215900675332435
$ws.Range("M9").Value = 0.4152066666666667
$ws.Range("N9").Value = 1.24562
$ws.Range("O9").Value = 0.02894309668502767
$ws.Range("P9").Value = 0.02894309668502767
$ws.Range("Q9").Value = 701.4041820567111
$ws.Range("R9").Value = 6312.637638510399
$ws.Range("S9").Value = 0.02667367042857585
$ws.Range("T9").Value = 0.02667367042857584
$ws.Range("G10").Value = 93.641553
$ws.Range("H10").Value = 280.924659
$ws.Range("I10").Value = 0.05108605424341119
$ws.Range("J10").Value = 0.05108605424341119
$ws.Range("M10").Value = 0.7489546666666667
$ws.Range("N10").Value = 2.246864
$ws.Range("O10").Value = 0.05220789806691288
$ws.Range("P10").Value = 0.05220789806691287
$ws.Range("Q10").Value = 70.133278113264
$ws.Range("R10").Value = 631.199503019376
$ws.Range("S10").Value = 0.002667095512580793
$ws.Range("T10").Value = 0.002667095512580793
$ws.Range("G11").Value = 93.641553
$ws.Range("H11").Value = 280.924659
$ws.Range("I11").Value = 0.05108605424341119
$ws.Range("J11").Value = 0.05108605424341119
$ws.Range("O11").Value = 0.1982273102638064
$ws.Range("P11").Value = 0.1982273102638064
$ws.Range("Q11").Value = 266.287891202931
$ws.Range("R11").Value = 2396.591020826379
$ws.Range("S11").Value = 0.01012665112466231
$ws.Range("T11").Value = 0.01012665112466231
$ws.Range("G12").Value = 93.641553
$ws.Range("H12").Value = 280.924659
$ws.Range("I12").Value = 0.05108605424341119
$ws.Range("J12").Value = 0.05108605424341119
$ws.Range("M12").Value = 10.337765
$ws.Range("N12").Value = 31.013295
$ws.Range("O12").Value = 0.7206216949842531
$ws.Range("P12").Value = 0.720621694984253
$ws.Range("Q12").Value = 968.044369149045
$ws.Range("R12").Value = 8712.399322341405
$ws.Range("S12").Value = 0.03681371899894446
$ws.Range("T12").Value = 0.03681371899894446
$ws.Range("G13").Value = 93.641553
$ws.Range("H13").Value = 280.924659
$ws.Range("I13").Value = 0.05108605424341119
$ws.Range("J13").Value = 0.05108605424341119
$ws.Range("M13").Value = 0.4152066666666667
$ws.Range("N13").Value = 1.24562
$ws.Range("O13").Value = 0.02894309668502767
$ws.Range("P13").Value = 0.02894309668502767
$ws.Range("Q13").Value = 38.88059708262
$ws.Range("R13").Value = 349.92537374358
$ws.Range("S13").Value = 0.001478588607223618
$ws.Range("T13").Value = 0.001478588607223618
$ws.Range("G14").Value = 14.34625366666667
$ws.Range("H14").Value = 43.038761
$ws.Range("I14").Value = 0.007826584134129748
$ws.Range("J14").Value = 0.007826584134129748
$ws.Range("M14").Value = 0.7489546666666667
$ws.Range("N14").Value = 2.246864
$ws.Range("O14").Value = 0.05220789806691288
$ws.Range("P14").Value = 0.05220789806691287
$ws.Range("Q14").Value = 10.74469363283378
$ws.Range("R14").Value = 96.702242695504
$ws.Range("S14").Value = 0.0004086095066867635
$ws.Range("T14").Value = 0.0004086095066867634
$ws.Range("G15").Value = 14.34625366666667
$ws.Range("H15").Value = 43.038761
$ws.Range("I15").Value = 0.007826584134129748
$ws.Range("J15").Value = 0.007826584134129748
$ws.Range("O15").Value = 0.1982273102638064
$ws.Range("P15").Value = 0.1982273102638064
$ws.Range("Q15").Value = 40.79635069229345
$ws.Range("R15").Value = 367.167156230641
$ws.Range("S15").Value = 0.001551442721461922
$ws.Range("T15").Value = 0.001551442721461922
$ws.Range("G16").Value = 14.34625366666667
$ws.Range("H16").Value = 43.038761
$ws.Range("I16").Value = 0.007826584134129748
$ws.Range("J16").Value = 0.007826584134129748
$ws.Range("M16").Value = 10.337765
$ws.Range("N16").Value = 31.013295
$ws.Range("O16").Value = 0.7206216949842531
$ws.Range("P16").Value = 0.720621694984253
$ws.Range("Q16").Value = 148.3081990363883
$ws.Range("R16").Value = 1334.773791327495
$ws.Range("S16").Value = 0.005640006324673442
$ws.Range("T16").Value = 0.005640006324673441
$ws.Range("G17").Value = 14.34625366666667
$ws.Range("H17").Value = 43.038761
$ws.Range("I17").Value = 0.007826584134129748
$ws.Range("J17").Value = 0.007826584134129748
$ws.Range("M17").Value = 0.4152066666666667
$ws.Range("N17").Value = 1.24562
$ws.Range("O17").Value = 0.02894309668502767
$ws.Range("P17").Value = 0.02894309668502767
$ws.Range("Q17").Value = 5.956660164091112
$ws.Range("R17").Value = 53.60994147682
$ws.Range("S17").Value = 0.0002265255813076209
$ws.Range("T17").Value = 0.0002265255813076208
